$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.517.58"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.848.22"
$s = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = $s
$ws.Range("E4").Value = "  +0.11%  "
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.54"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  +0.03%  "
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6298"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07452"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = "  -1.79%  "
$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2906"
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = "  -0.43%  "
$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.97"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  +1.81%  "
$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07743"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "1.843.77"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  +0.00%  "
$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6823"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = "  +0.57%  "
$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001023"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = "  -1.97%  "
$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.58"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  -0.72%  "
$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.319"
$ws.Range("D17").Style = $s
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("D18").Value = "29.534.74"
$ws.Range("E18").Value = "  +0.54%  "
$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.74"
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = "  -0.15%  "
$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.36"
$ws.Range("D20").Style = $s
$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  -0.01%  "
$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.513"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").Value = "  +0.09%  "
$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.29"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +0.80%  "
$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1362"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("E27").Value = "  -0.71%  "
$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06591"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  +15.92%  "
$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.462"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  +2.55%  "
$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.488"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  +0.96%  "
$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.097"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  -0.50%  "
$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.089"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  +1.13%  "
$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.847"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  -1.24%  "
$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6971"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -0.58%  "
$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01872"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  +2.56%  "
$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.839"
$ws.Range("D38").Style = $s
$ws.Range("D39").Value = "1.252.90"
$ws.Range("E39").Value = "  +1.23%  "
$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.788"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  +5.63%  "
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9383"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "2.004.89"
$ws.Range("E43").Value = "  -0.16%  "
$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.22"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  -0.23%  "
$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.22"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  +0.75%  "
$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.095"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  +2.73%  "
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1157"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  -0.08%  "
$s = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.017"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = "  -0.37%  "
$s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3934"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = "  -1.25%  "
$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000112"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -4.75%  "

Write-Host "Applied changes"